$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03247682266724
$ws.Cells.Item(2, 4).Value = 1.041880590087326
$ws.Cells.Item(2, 5).Value = 1.03186514798708
$ws.Cells.Item(2, 6).Value = 1.049050770049192
$ws.Cells.Item(2, 9).Value = 1.027160259577969
$ws.Cells.Item(2, 10).Value = 1.037606830111397
$ws.Cells.Item(2, 11).Value = 1.04465885329997
$ws.Cells.Item(2, 12).Value = 1.034672032303414
$ws.Cells.Item(2, 13).Value = 1.051808907151408
$ws.Cells.Item(2, 14).Value = 1.039080350825322

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033907953514061
$ws.Cells.Item(3, 4).Value = 1.043245750382815
$ws.Cells.Item(3, 5).Value = 1.033096635116719
$ws.Cells.Item(3, 6).Value = 1.050463900759419
$ws.Cells.Item(3, 9).Value = 1.027151154978046
$ws.Cells.Item(3, 10).Value = 1.038677826411654
$ws.Cells.Item(3, 11).Value = 1.04583327782716
$ws.Cells.Item(3, 12).Value = 1.035711040097946
$ws.Cells.Item(3, 13).Value = 1.053032653776008
$ws.Cells.Item(3, 14).Value = 1.04015286806317

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034833265293881
$ws.Cells.Item(4, 4).Value = 1.044128490390494
$ws.Cells.Item(4, 5).Value = 1.03389310669884
$ws.Cells.Item(4, 6).Value = 1.051377070553262
$ws.Cells.Item(4, 9).Value = 1.027142772198315
$ws.Cells.Item(4, 10).Value = 1.039369718705235
$ws.Cells.Item(4, 11).Value = 1.046592071287426
$ws.Cells.Item(4, 12).Value = 1.03638241121655
$ws.Cells.Item(4, 13).Value = 1.053822760243242
$ws.Cells.Item(4, 14).Value = 1.040845742923169

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035222098228875
$ws.Cells.Item(5, 4).Value = 1.044499452715199
$ws.Cells.Item(5, 5).Value = 1.034227855523841
$ws.Cells.Item(5, 6).Value = 1.051760679387014
$ws.Cells.Item(5, 9).Value = 1.027138650937356
$ws.Cells.Item(5, 10).Value = 1.039660327327636
$ws.Cells.Item(5, 11).Value = 1.046910799642383
$ws.Cells.Item(5, 12).Value = 1.036664434974856
$ws.Cells.Item(5, 13).Value = 1.054154508375991
$ws.Cells.Item(5, 14).Value = 1.041136764243147

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035287375291288
$ws.Cells.Item(6, 4).Value = 1.044561730748932
$ws.Cells.Item(6, 5).Value = 1.034284056285683
$ws.Cells.Item(6, 6).Value = 1.051825072251442
$ws.Cells.Item(6, 9).Value = 1.027137923923094
$ws.Cells.Item(6, 10).Value = 1.039709106487042
$ws.Cells.Item(6, 11).Value = 1.046964299910309
$ws.Cells.Item(6, 12).Value = 1.036711775167052
$ws.Cells.Item(6, 13).Value = 1.054210186260434
$ws.Cells.Item(6, 14).Value = 1.041185612674555

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034838461552381
$ws.Cells.Item(7, 4).Value = 1.044133447758046
$ws.Cells.Item(7, 5).Value = 1.033897579969389
$ws.Cells.Item(7, 6).Value = 1.051382197478574
$ws.Cells.Item(7, 9).Value = 1.027142719477056
$ws.Cells.Item(7, 10).Value = 1.039373602860239
$ws.Cells.Item(7, 11).Value = 1.04659633120114
$ws.Cells.Item(7, 12).Value = 1.036386180494534
$ws.Cells.Item(7, 13).Value = 1.053827194694758
$ws.Cells.Item(7, 14).Value = 1.040849632594119

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.032960633630022
$ws.Cells.Item(8, 4).Value = 1.042342081393377
$ws.Cells.Item(8, 5).Value = 1.032281416440214
$ws.Cells.Item(8, 6).Value = 1.049528599502753
$ws.Cells.Item(8, 9).Value = 1.027157698117068
$ws.Cells.Item(8, 10).Value = 1.037969011252064
$ws.Cells.Item(8, 11).Value = 1.045055993535343
$ws.Cells.Item(8, 12).Value = 1.035023365862268
$ws.Cells.Item(8, 13).Value = 1.052222840506162
$ws.Cells.Item(8, 14).Value = 1.039443046304752

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029645853915091
$ws.Cells.Item(9, 4).Value = 1.039180592639943
$ws.Cells.Item(9, 5).Value = 1.029430423845227
$ws.Cells.Item(9, 6).Value = 1.046252779519013
$ws.Cells.Item(9, 9).Value = 1.027165039368521
$ws.Cells.Item(9, 10).Value = 1.035485240010899
$ws.Cells.Item(9, 11).Value = 1.042332829705012
$ws.Cells.Item(9, 12).Value = 1.032614582419648
$ws.Cells.Item(9, 13).Value = 1.049382270193302
$ws.Cells.Item(9, 14).Value = 1.036955747823532

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027431745170299
$ws.Cells.Item(10, 4).Value = 1.037069374622477
$ws.Cells.Item(10, 5).Value = 1.027527429492147
$ws.Cells.Item(10, 6).Value = 1.044062203185137
$ws.Cells.Item(10, 9).Value = 1.027157160870885
$ws.Cells.Item(10, 10).Value = 1.033823302360875
$ws.Cells.Item(10, 11).Value = 1.040511156793348
$ws.Cells.Item(10, 12).Value = 1.031003575219406
$ws.Cells.Item(10, 13).Value = 1.047479241620004
$ws.Cells.Item(10, 14).Value = 1.035291450031417

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02647192355107
$ws.Cells.Item(11, 4).Value = 1.036154286070151
$ws.Cells.Item(11, 5).Value = 1.026702803597876
$ws.Cells.Item(11, 6).Value = 1.043112007877939
$ws.Cells.Item(11, 9).Value = 1.027150727360527
$ws.Cells.Item(11, 10).Value = 1.033102166486809
$ws.Cells.Item(11, 11).Value = 1.039720817755571
$ws.Cells.Item(11, 12).Value = 1.030304719098081
$ws.Cells.Item(11, 13).Value = 1.046652944738684
$ws.Cells.Item(11, 14).Value = 1.034569290061695

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026115231257817
$ws.Cells.Item(12, 4).Value = 1.035814238051431
$ws.Cells.Item(12, 5).Value = 1.026396403333128
$ws.Cells.Item(12, 6).Value = 1.042758807865801
$ws.Cells.Item(12, 9).Value = 1.027147884175114
$ws.Cells.Item(12, 10).Value = 1.032834073530135
$ws.Cells.Item(12, 11).Value = 1.039427013892046
$ws.Cells.Item(12, 12).Value = 1.030044936188399
$ws.Cells.Item(12, 13).Value = 1.046345674924673
$ws.Cells.Item(12, 14).Value = 1.034300816382272

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02619175085076
$ws.Cells.Item(13, 4).Value = 1.035887186084366
$ws.Cells.Item(13, 5).Value = 1.026462131719553
$ws.Cells.Item(13, 6).Value = 1.042834582104234
$ws.Cells.Item(13, 9).Value = 1.027148514561073
$ws.Cells.Item(13, 10).Value = 1.032891590908382
$ws.Cells.Item(13, 11).Value = 1.039490046612278
$ws.Cells.Item(13, 12).Value = 1.030100669473556
$ws.Cells.Item(13, 13).Value = 1.046411601095131
$ws.Cells.Item(13, 14).Value = 1.034358415441794

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.026442442799488
$ws.Cells.Item(14, 4).Value = 1.036126180552325
$ws.Cells.Item(14, 5).Value = 1.026677478473243
$ws.Cells.Item(14, 6).Value = 1.043082817492516
$ws.Cells.Item(14, 9).Value = 1.027150501589999
$ws.Cells.Item(14, 10).Value = 1.033080010578928
$ws.Cells.Item(14, 11).Value = 1.039696536695289
$ws.Cells.Item(14, 12).Value = 1.030283249398585
$ws.Cells.Item(14, 13).Value = 1.046627552814032
$ws.Cells.Item(14, 14).Value = 1.034547102689883

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026596879410901
$ws.Cells.Item(15, 4).Value = 1.03627341376406
$ws.Cells.Item(15, 5).Value = 1.02681014765203
$ws.Cells.Item(15, 6).Value = 1.043235729528268
$ws.Cells.Item(15, 9).Value = 1.027151665787691
$ws.Cells.Item(15, 10).Value = 1.033196071406546
$ws.Cells.Item(15, 11).Value = 1.039823730528215
$ws.Cells.Item(15, 12).Value = 1.030395716728239
$ws.Cells.Item(15, 13).Value = 1.04676056175251
$ws.Cells.Item(15, 14).Value = 1.03466332833719

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027495421653768
$ws.Cells.Item(16, 4).Value = 1.03713008623543
$ws.Cells.Item(16, 5).Value = 1.027582143790824
$ws.Cells.Item(16, 6).Value = 1.044125229041678
$ws.Cells.Item(16, 9).Value = 1.027157524216574
$ws.Cells.Item(16, 10).Value = 1.033871129642387
$ws.Cells.Item(16, 11).Value = 1.040563576065506
$ws.Cells.Item(16, 12).Value = 1.031049928667221
$ws.Cells.Item(16, 13).Value = 1.047534031918939
$ws.Cells.Item(16, 14).Value = 1.035339345233155

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028058754690353
$ws.Cells.Item(17, 4).Value = 1.037667204747761
$ws.Cells.Item(17, 5).Value = 1.028066228045587
$ws.Cells.Item(17, 6).Value = 1.044682739886002
$ws.Cells.Item(17, 9).Value = 1.027160390345438
$ws.Cells.Item(17, 10).Value = 1.034294169329722
$ws.Cells.Item(17, 11).Value = 1.041027245102741
$ws.Cells.Item(17, 12).Value = 1.031459952898278
$ws.Cells.Item(17, 13).Value = 1.048018597355958
$ws.Cells.Item(17, 14).Value = 1.03576298568537

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028387231852338
$ws.Cells.Item(18, 4).Value = 1.037980408770541
$ws.Cells.Item(18, 5).Value = 1.028348527019498
$ws.Cells.Item(18, 6).Value = 1.045007766617997
$ws.Cells.Item(18, 9).Value = 1.02716177048917
$ws.Cells.Item(18, 10).Value = 1.03454077624674
$ws.Cells.Item(18, 11).Value = 1.041297546904369
$ws.Cells.Item(18, 12).Value = 1.031698989955104
$ws.Cells.Item(18, 13).Value = 1.048301017202988
$ws.Cells.Item(18, 14).Value = 1.036009942812496

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028499216393339
$ws.Cells.Item(19, 4).Value = 1.03808718846997
$ws.Cells.Item(19, 5).Value = 1.028444773781221
$ws.Cells.Item(19, 6).Value = 1.045118565377996
$ws.Cells.Item(19, 9).Value = 1.027162191601837
$ws.Cells.Item(19, 10).Value = 1.034624838384764
$ws.Cells.Item(19, 11).Value = 1.041389687794825
$ws.Cells.Item(19, 12).Value = 1.03178047466039
$ws.Cells.Item(19, 13).Value = 1.048397278077052
$ws.Cells.Item(19, 14).Value = 1.036094124328398

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027998325345644
$ws.Cells.Item(20, 4).Value = 1.037609586169688
$ws.Cells.Item(20, 5).Value = 1.028014296542926
$ws.Cells.Item(20, 6).Value = 1.044622940829542
$ws.Cells.Item(20, 9).Value = 1.027160112996758
$ws.Cells.Item(20, 10).Value = 1.034248796204799
$ws.Cells.Item(20, 11).Value = 1.04097751319981
$ws.Cells.Item(20, 12).Value = 1.031415973947788
$ws.Cells.Item(20, 13).Value = 1.047966630755436
$ws.Cells.Item(20, 14).Value = 1.035717548125405

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026368625066941
$ws.Cells.Item(21, 4).Value = 1.036055806641493
$ws.Cells.Item(21, 5).Value = 1.026614066951824
$ws.Cells.Item(21, 6).Value = 1.043009725446648
$ws.Cells.Item(21, 9).Value = 1.02714992897389
$ws.Cells.Item(21, 10).Value = 1.033024532096718
$ws.Cells.Item(21, 11).Value = 1.039635737088326
$ws.Cells.Item(21, 12).Value = 1.03022948962493
$ws.Cells.Item(21, 13).Value = 1.046563970015906
$ws.Cells.Item(21, 14).Value = 1.034491545421865

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025342970364464
$ws.Cells.Item(22, 4).Value = 1.035078051400211
$ws.Cells.Item(22, 5).Value = 1.02573312056361
$ws.Cells.Item(22, 6).Value = 1.041993953791423
$ws.Cells.Item(22, 9).Value = 1.027140902358004
$ws.Cells.Item(22, 10).Value = 1.032253449884341
$ws.Cells.Item(22, 11).Value = 1.038790736902033
$ws.Cells.Item(22, 12).Value = 1.029482359989818
$ws.Cells.Item(22, 13).Value = 1.045680055151551
$ws.Cells.Item(22, 14).Value = 1.033719368184302

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025886786175173
$ws.Cells.Item(23, 4).Value = 1.035596458727541
$ws.Cells.Item(23, 5).Value = 1.026200182226417
$ws.Cells.Item(23, 6).Value = 1.042532575640399
$ws.Cells.Item(23, 9).Value = 1.027145936040995
$ws.Cells.Item(23, 10).Value = 1.032662343645356
$ws.Cells.Item(23, 11).Value = 1.039238819240733
$ws.Cells.Item(23, 12).Value = 1.029878537072805
$ws.Cells.Item(23, 13).Value = 1.046148826965547
$ws.Cells.Item(23, 14).Value = 1.034128842621358

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028025631091584
$ws.Cells.Item(24, 4).Value = 1.037635621797619
$ws.Cells.Item(24, 5).Value = 1.028037762336836
$ws.Cells.Item(24, 6).Value = 1.04464996194369
$ws.Cells.Item(24, 9).Value = 1.027160239219924
$ws.Cells.Item(24, 10).Value = 1.034269298815154
$ws.Cells.Item(24, 11).Value = 1.040999985365267
$ws.Cells.Item(24, 12).Value = 1.031435846523257
$ws.Cells.Item(24, 13).Value = 1.047990112902358
$ws.Cells.Item(24, 14).Value = 1.035738079851818

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030503528130116
$ws.Cells.Item(25, 4).Value = 1.039998519673605
$ws.Cells.Item(25, 5).Value = 1.030167866896008
$ws.Cells.Item(25, 6).Value = 1.047100817133351
$ws.Cells.Item(25, 9).Value = 1.02716539543726
$ws.Cells.Item(25, 10).Value = 1.03612841016152
$ws.Cells.Item(25, 11).Value = 1.043037912400111
$ws.Cells.Item(25, 12).Value = 1.033238202211382
$ws.Cells.Item(25, 13).Value = 1.05011824960771
$ws.Cells.Item(25, 14).Value = 1.037599831349539
